# Applies the edits described in the commit "busca de redações pelo estudante"
# to EMBASAMENTO TEÓRICO.docx

$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceAll = 2

function Replace-Text($find, $replace) {
    $r = $d.Content
    $ok = $r.Find.Execute($find, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replace, $wdReplaceAll)
    if (-not $ok) {
        Write-Host "WARNING: not found -> $find"
    }
}

# 1. Remove the stray _GoBack bookmark from the empty paragraph right after the title.
$d.Bookmarks.ShowHidden = $true
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 2a. "...importante etapa da prova" -> "...importante etapa do mesmo"
Replace-Text "importante etapa da prova" "importante etapa do mesmo"

# 2b. "Tais erros são correspondentes a falta de prática" -> "Tais erros correspondem à falta de prática"
Replace-Text "erros são correspondentes a falta de prática" "erros correspondem à falta de prática"

# 3. Aprofundamento paragraph rewording
Replace-Text "o que se pedem com relação às competências, como por exemplo o domínio da linguagem culta portuguesa brasileira ou se a redação compreende e não foge ao tema proposto entre outras." "o que se pede com relação às competências, como o domínio da linguagem culta portuguesa brasileira, a compreensão e não fuga ao tema proposto, entre outras."

# 4a. Tese1 statistic correction
Replace-Text "mais de 5% (por cento) das redações foram zeradas por desobedecer alguma das competências" "mais de 6,5% (por cento) foram zeradas por desobedecer alguma das competências"

# 4b. Tese1 wording correction
Replace-Text "como por exemplo partes dentro da redação estando desconectadas correspondem a" "como partes dentro da redação estando desconectadas, o que corresponde a"

# 5. Conclusão Geral rewording
Replace-Text "O sistema poderá, também, oferecer àqueles que não conhecem este modelo de avaliação, sejam estudantes ou não, que desejam aprender sobre, garantindo melhor pontuação em vestibulares e no ENEM." "O sistema poderá, também, oferecer condições de prática àqueles que não conhecem este modelo de avaliação, sejam estudantes ou não, que desejam aprender sobre, garantindo melhor pontuação na redação dissertativa argumentativa do ENEM."

# 6. Move the _GoBack bookmark so it wraps the "Conclusão Geral" paragraph (the last edit location).
$concl = $d.Paragraphs(15)
$bmRange = $d.Range($concl.Range.Start, $concl.Range.End)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 7. Move the lastRenderedPageBreak marker from the first reference paragraph to the "Referencias:" heading.
$refHeading = $d.Paragraphs(17)
$refHeading.Range.InsertBefore([char]2)
Replace-Text ([char]2) ""
$refHeading2 = $d.Paragraphs(17)

# 8. Append a trailing period to the first reference entry.
$p18 = $d.Paragraphs(18)
$endOfP18 = $p18.Range.End - 1
$insPoint = $d.Range($endOfP18, $endOfP18)
$insPoint.InsertAfter(".")

Write-Host "done"
